# The edit only touches the <w:nsid> values inside the document's
# numbering definitions (word/numbering.xml) -- six abstractNum entries
# get a freshly minted nsid GUID. Nothing else in the package changes
# (no visible text/formatting edits), so this is applied by round-
# tripping the document's OOXML through Range.XML()/InsertXML() with
# the six old nsid values swapped for their new ones.

$d = $word.ActiveDocument

# old nsid -> new nsid, one pair per <w:abstractNum> definition touched
# by the change (abstractNumId 990, 99411, 991, 99417, 994113, 99416)
$nsidMap = @(
    @("2c594fbb", "195d14c0"),   # abstractNumId 990
    @("b367231d", "e823aec2"),   # abstractNumId 99411
    @("379fdca2", "3ee8819c"),   # abstractNumId 991
    @("858d8ee4", "b77d9be5"),   # abstractNumId 99417
    @("1ac62ab9", "5c6f3503"),   # abstractNumId 994113
    @("b32afb94", "d8548c3a")    # abstractNumId 99416
)

$range = $d.Content
$xml = $range.XML()

foreach ($pair in $nsidMap) {
    $oldNsid = $pair[0]
    $newNsid = $pair[1]
    $needle = '<w:nsid w:val="' + $oldNsid + '"'
    $replacement = '<w:nsid w:val="' + $newNsid + '"'
    if ($xml.IndexOf($needle) -lt 0) {
        throw "nsid pattern not found: $needle"
    }
    $xml = $xml.Replace($needle, $replacement)
}

$range.InsertXML($xml)

Write-Output "nsid values updated"
